# Updated cryptos list on Fri Oct 18 23:40:31 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain stored as text (the Price column
# holds strings like "599.96" / "68.360.00" that Excel would otherwise
# auto-coerce into numbers). Force text format, assign, then restore the
# plain "Normal" style so no stray number-format styling is left behind.
function Set-Text($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Pct($row, $value) {
    $ws.Range("E$row").Value = "  $value  "
}

# Row 2 - Bitcoin
Set-Text "D2" "68.360.00"
Set-Pct 2 "+1.50%"

# Row 3 - Ethereum
Set-Text "D3" "2.641.33"
Set-Pct 3 "+1.44%"

# Row 4 - TetherUSD
Set-Pct 4 "-0.01%"

# Row 5 - BNB
Set-Text "D5" "599.96"
Set-Pct 5 "+1.26%"

# Row 6 - Solana
Set-Text "D6" "154.78"
Set-Pct 6 "+2.94%"

# Row 7 - USDC
Set-Pct 7 "-0.01%"

# Row 8 - XRP
Set-Pct 8 "+0.45%"

# Row 9 - LidoStakedEther
Set-Text "D9" "2.640.63"
Set-Pct 9 "+1.43%"

# Row 10 - Dogecoin
Set-Pct 10 "+6.77%"

# Row 11 - TRON
Set-Pct 11 "-0.56%"

# Row 12 - Toncoin
Set-Pct 12 "+1.36%"

# Row 13 - Cardano
Set-Pct 13 "+2.02%"

# Row 14 - Avalanche
Set-Pct 14 "+3.04%"

# Row 15 - ShibaInu
Set-Pct 15 "+3.40%"

# Row 16 - WrappedliquidstakedEther2.0
Set-Text "D16" "3.120.85"
Set-Pct 16 "+1.37%"

# Row 17 - WrappedBTC
Set-Text "D17" "68.288.63"

# Row 18 - WrappedEther
Set-Text "D18" "2.641.62"
Set-Pct 18 "+1.38%"

# Row 19 - Chainlink
Set-Pct 19 "+4.04%"

# Row 20 - BitcoinCash
Set-Text "D20" "366.83"
Set-Pct 20 "-0.90%"

# Row 21 - Uniswap
Set-Text "D21" "7.45"
Set-Pct 21 "+1.47%"

# Row 22 - Polkadot
Set-Text "D22" "4.30"
Set-Pct 22 "+2.68%"

# Row 23 - NEARProtocol
Set-Text "D23" "4.87"
Set-Pct 23 "+2.47%"

# Row 24 - SuiNetwork
Set-Pct 24 "+4.94%"

# Row 25 - Litecoin
Set-Text "D25" "73.51"
Set-Pct 25 "+0.46%"

# Row 27 - Aptos
Set-Pct 27 "+1.30%"

# Row 28 - WrappedeETH
Set-Text "D28" "2.772.69"
Set-Pct 28 "+1.48%"

# Row 29 - PEPE
Set-Pct 29 "+6.31%"

# Row 30 - Binance-PegBSC-USD
Set-Text "D30" "0.999"
Set-Pct 30 "-0.08%"

# Row 31 - Bittensor
Set-Text "D31" "574.23"
Set-Pct 31 "-0.53%"

# Row 32 - Fetch.AI
Set-Pct 32 "+5.11%"

# Row 33 - InternetComputer(DFINITY)
Set-Pct 33 "+4.73%"

# Row 34 - PancakeSwap
Set-Pct 34 "+2.75%"

# Row 36 - FirstDigitalUSD
Set-Text "D36" "0.999"

# Row 37 - ImmutableX
Set-Pct 37 "+3.70%"

# Row 38 - Monero
Set-Text "D38" "160.82"
Set-Pct 38 "+1.53%"

# Row 39 - EthereumClassic
Set-Pct 39 "+1.59%"

# Row 40 - Stacks
Set-Text "D40" "1.92"
Set-Pct 40 "+3.37%"

# Row 41 - PolygonEcosystemToken
Set-Pct 41 "+1.21%"

# Row 42 - RenderToken
Set-Pct 42 "+3.86%"

# Row 43 - dogwifhat
Set-Pct 43 "+4.03%"

# Row 44 - WhiteBITCoin
Set-Pct 44 "+3.63%"

# Row 45 - BabyDogeCoin
Set-Text "D45" "0.0₆0320"
Set-Pct 45 "+14.00%"

# Row 46 - USDe
Set-Pct 46 "+0.05%"

# Row 47 - OKB
Set-Pct 47 "-0.26%"

# Row 48 - Aave
Set-Text "D48" "158.44"
Set-Pct 48 "+3.69%"

# Row 49 - Filecoin
Set-Pct 49 "+3.53%"

# Row 50 - Optimism
Set-Text "D50" "1.72"
Set-Pct 50 "+2.66%"

# Row 51 - InjectiveProtocol
Set-Text "D51" "22.02"
Set-Pct 51 "+3.49%"
